$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new user rows first (so their new strings occupy the shared
# string table slots ahead of the new "rol" header string)
$ws.Range("A7").Value = "admin_1vg"
$ws.Range("B7").Value = 1234
$ws.Range("C7").Value = "Imelda Garza"
$ws.Range("D7").Value = "1VG"
$ws.Range("E7").Value = "usuario"

$ws.Range("A8").Value = "admin_3vg"
$ws.Range("B8").Value = 1234
$ws.Range("C8").Value = "Angelica Cruz"
$ws.Range("D8").Value = "3VG"
$ws.Range("E8").Value = "usuario"

# Fill in "rol" values for existing rows
$ws.Range("E2").Value = "admin"
$ws.Range("E3").Value = "usuario"
$ws.Range("E4").Value = "usuario"
$ws.Range("E5").Value = "usuario"
$ws.Range("E6").Value = "admin"

# Add new "rol" column header last
$ws.Range("E1").Value = "rol"

# Update selection to match target state
$ws.Range("E9").Select()
